$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.273.97"
$ws.Range("E2").Value = "'  +0.37%  "
$ws.Range("D3").Value = "'1.901.59"
$ws.Range("E3").Value = "'  -0.37%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.21%  "
$ws.Range("D5").Value = "'326.48"
$ws.Range("E5").Value = "'  -0.33%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D7").Value = "'0.4651"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("D8").Value = "'0.3929"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("D9").Value = "'0.07892"
$ws.Range("E9").Value = "'  -0.88%  "
$ws.Range("D10").Value = "'0.9890"
$ws.Range("E10").Value = "'  -1.38%  "
$ws.Range("D11").Value = "'22.03"
$ws.Range("E11").Value = "'  -1.45%  "
$ws.Range("D12").Value = "'1.924.18"
$ws.Range("E12").Value = "'  +0.18%  "
$ws.Range("D13").Value = "'7.080"
$ws.Range("E13").Value = "'  -0.81%  "
$ws.Range("D14").Value = "'5.752"
$ws.Range("E14").Value = "'  -0.83%  "
$ws.Range("D15").Value = "'0.06989"
$ws.Range("E15").Value = "'  +0.46%  "
$ws.Range("D16").Value = "'88.41"
$ws.Range("E16").Value = "'  -0.46%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "'  -0.01%  "
$ws.Range("D18").Value = "'0.000009977"
$ws.Range("E18").Value = "'  -1.29%  "
$ws.Range("D19").Value = "'17.09"
$ws.Range("E19").Value = "'  -0.73%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "'  -0.15%  "
$ws.Range("D21").Value = "'29.279.55"
$ws.Range("E21").Value = "'  +0.29%  "
$ws.Range("D22").Value = "'5.319"
$ws.Range("E22").Value = "'  -0.97%  "
$ws.Range("D23").Value = "'11.09"
$ws.Range("E23").Value = "'  +0.15%  "
$ws.Range("D24").Value = "'2.095"
$ws.Range("E24").Value = "'  +1.78%  "
$ws.Range("D25").Value = "'156.22"
$ws.Range("D26").Value = "'19.47"
$ws.Range("E26").Value = "'  -0.25%  "
$ws.Range("D27").Value = "'5.983"
$ws.Range("D28").Value = "'118.59"
$ws.Range("E28").Value = "'  -0.75%  "
$ws.Range("D29").Value = "'1.908"
$ws.Range("E29").Value = "'  -5.00%  "
$ws.Range("D30").Value = "'0.09373"
$ws.Range("D31").Value = "'0.9083"
$ws.Range("E31").Value = "'  -1.73%  "
$ws.Range("D32").Value = "'5.287"
$ws.Range("E32").Value = "'  -1.44%  "
$ws.Range("D33").Value = "'1.329"
$ws.Range("E33").Value = "'  -1.35%  "
$ws.Range("D34").Value = "'3.213"
$ws.Range("E34").Value = "'  -1.56%  "
$ws.Range("D35").Value = "'1.189"
$ws.Range("E35").Value = "'  +2.05%  "
$ws.Range("D36").Value = "'0.05780"
$ws.Range("E36").Value = "'  -1.27%  "
$ws.Range("D37").Value = "'0.02090"
$ws.Range("E37").Value = "'  -0.67%  "
$ws.Range("D38").Value = "'1.000"
$ws.Range("E38").Value = "'  -0.17%  "
$ws.Range("D39").Value = "'7.743"
$ws.Range("E39").Value = "'  -3.63%  "
$ws.Range("D40").Value = "'0.5712"
$ws.Range("E40").Value = "'  -0.93%  "
$ws.Range("D41").Value = "'0.1787"
$ws.Range("E41").Value = "'  -1.45%  "
$ws.Range("D42").Value = "'9.771"
$ws.Range("E42").Value = "'  -2.55%  "
$ws.Range("D43").Value = "'11.90"
$ws.Range("E43").Value = "'  -1.59%  "
$ws.Range("D44").Value = "'0.5358"
$ws.Range("E44").Value = "'  -1.36%  "
$ws.Range("D45").Value = "'2.193"
$ws.Range("E45").Value = "'  -1.73%  "
$ws.Range("D46").Value = "'0.07039"
$ws.Range("E46").Value = "'  -0.82%  "
$ws.Range("D47").Value = "'1.854"
$ws.Range("E47").Value = "'  -1.59%  "
$ws.Range("E48").Value = "'  +0.45%  "
$ws.Range("D49").Value = "'113.56"
$ws.Range("E49").Value = "'  +0.86%  "
$ws.Range("E50").Value = "'  -2.53%  "
$ws.Range("E51").Value = "'  -0.44%  "
